$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PurchaseList")

# --- Designator / Comment / Footprint edits (board rerouted: connector & resistor renumbering) ---
# Leading "'" forces the literal-text quote-prefix cell style to stick (matches
# the original cells, which were all entered as quoted text) instead of Excel
# reassigning a fresh style without quotePrefix.

# Row 4: JP1, JP2 -> J2, J3
$ws.Range("B4").Value = "'J2, J3"

# Row 5: R1 -> R1, R3, R4, R5 (added pull-ups), Quantity 1 -> 4
$ws.Range("B5").Value = "'R1, R3, R4, R5"
$ws.Range("G5").Value = 4

# Row 9 and Row 10 swap content: JP3/HEADER_5/1X05 <-> LCD1/LCD MODULE SIL/1X14-FEMALE
$ws.Range("B9").Value = "'J1"
$ws.Range("C9").Value = "'LCD MODULE SIL"
$ws.Range("D9").Value = "'1X14-FEMALE"

$ws.Range("B10").Value = "'J4"
$ws.Range("C10").Value = "'HEADER_5"
$ws.Range("D10").Value = "'1X05"

# --- Quotation timestamp (new auto quotation) ---
$ws.Range("E14").Value = "'21:44"
